$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column D values (rows 2-9)
$ws.Range("D2").Value = 89.51015903888973
$ws.Range("D3").Value = 94.14175759316585
$ws.Range("D4").Value = 92.30137333371246
$ws.Range("D5").Value = 85.97876663481171
$ws.Range("D6").Value = 82.26939119452892
$ws.Range("D7").Value = 82.08288518084234
$ws.Range("D8").Value = 76.1946382788903
$ws.Range("D9").Value = 70.97284275433343

# Update column C values (rows 6-9)
$ws.Range("C6").Value = 2.941176470588236
$ws.Range("C7").Value = 2.083333333333333
$ws.Range("C8").Value = 0.9803921568627446
$ws.Range("C9").Value = 0.2777777777777774
